# Add 5 new "RSNA" event rows (rows 3-7) to the events worksheet, each a
# copy of the existing row 2 (NIAID BioVisualization / RSNA 2023 exhibit)
# with an updated title and its own hyperlink to the same RSNA URL.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url = "https://bioinformatics.niaid.nih.gov/rsna-2023"

$titles = @(
    "RSNA-2 2023",
    "RSNA -3 2023",
    "RSNA- 4 2023",
    "RSNA -5 2023",
    "RSNA - 6 2023"
)

# Add the hyperlinks to the destination cells first. Excel stamps the
# built-in "Hyperlink" cell style on the target when a hyperlink is
# inserted; doing this before the row copy means the subsequent copy
# (which carries row 2's own hyperlink-styled G cell) overwrites that
# stamp with the correct, already-in-use formatting.
for ($i = 0; $i -lt $titles.Length; $i++) {
    $row = 3 + $i
    $ws.Hyperlinks.Add($ws.Range("G" + $row), $url)
}

# Copy row 2's values/formats (A:K) down into the five new rows, matching
# row 2's own row height (wrapped description text makes the row tall).
for ($i = 0; $i -lt $titles.Length; $i++) {
    $row = 3 + $i
    $ws.Range("A2:K2").Copy($ws.Range("A" + $row + ":K" + $row))
    $ws.Rows($row).RowHeight = $ws.Rows(2).RowHeight
}

# Give each new row its own title in column B.
for ($i = 0; $i -lt $titles.Length; $i++) {
    $row = 3 + $i
    $ws.Range("B" + $row).Value2 = $titles[$i]
}

# Row 2's description cell (J2) still carries the older, font-duplicating
# style; nudging it through a self-copy collapses it onto the same
# wrap-text style already used by the rest of column J (J3:J7).
$ws.Range("J2").Copy($ws.Range("J2"))

# Match the saved selection left behind by the edit.
$ws.Range("B6").Select()
